$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where columns J (10) through AS (45) are reset to a constant 1
$onesRows = @(4,5,6,7,8,9,10,12,13,14,15,21,22)
foreach ($r in $onesRows) {
    for ($c = 10; $c -le 45; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

$row96 = @(4497955.98621076,6349316.81667575,6444216.92362503,5259717.40745011,5167069.16776774,4541190.79984833,4541191,4854451.97841785,5189322.3629583,5547292.81624818,5929956.82998119,6339017.81828481,6776296.70040275,7243740.00644687,7743428.54230102,8277586.65224837,8848592.120554339,9458986.756081739,10111487.7070551,10808999.5563414,11554627.2510887,12351689.9242794,13203735.6697224,14114557.3362561,15088209.4114695,16129026.0700972,17241640.4674316,18431005.3636354,19702415.1707624,21061529.5206296,22514398.4584498,24067489.3743755,25727715.7928366,27502468.147829,29399646.6811465,31427696.6100033)
for ($i = 0; $i -lt $row96.Length; $i++) {
    $ws.Cells.Item(96, 10 + $i).Value = $row96[$i]
}

$row97 = @(168315.742466078,161357.051248769,177195.393132506,163841.954758891,165493.036178719,177589.058037976,177589.1,184731.619678592,192161.406916731,199890.01543138,207929.463622585,216292.253262734,224991.388937481,234040.39826857,243453.352950017,253244.890630341,263430.237674895,274025.232843682,285046.351921468,296510.733338515,308436.204821755,320841.31111786,333745.342831319,347168.366422358,361131.255411358,375655.722838302,390764.355027712,406480.646711593,422829.037564999,439834.950211035,457524.829754402,475926.184904946,495067.630755175,514978.933278263,535691.055615729,557236.2062267819)
for ($i = 0; $i -lt $row97.Length; $i++) {
    $ws.Cells.Item(97, 10 + $i).Value = $row97[$i]
}

$row98 = @(374272.456707285,323395.116624136,274861.383079098,252767.436903689,252721.226218507,341278.456469139,341278.5,338102.926081881,334956.900669482,331840.148815908,328752.398132626,325693.37876566,322662.823372009,319660.467096279,316686.047547538,313739.304776388,310819.981252238,307927.821840804,305062.573781812,302223.986666902,299411.81241775,296625.805264382,293865.7217237,291131.320578199,288422.362854886,285738.611804397,283079.832880305,280445.793718623,277836.264117494,275251.016017074,272689.823479602,270152.462669653,267638.711834574,265148.351285108,262681.163376192,260236.932487933)
for ($i = 0; $i -lt $row98.Length; $i++) {
    $ws.Cells.Item(98, 10 + $i).Value = $row98[$i]
}

$row99 = @(30235.2939796777,28050.7280600012,31708.0748548513,36413.1228462395,36638.8406999203,27148.527354939,27148.53,27489.5996320865,27834.9541552494,28184.6474009933,28538.7338771135,28897.2687761921,29260.3079842009,29627.9080892124,30000.126390221,30377.0209060737,30758.650384514,31145.0743113389,31536.3529196715,31932.5471993494,32333.7189064314,32739.9305728239,33151.2455160275,33567.7278490069,33989.4424901842,34416.455173558,34848.8324589494,35286.6417423772,35729.951266563,36178.830131568,36633.3483055647,37093.5766357423,37559.5868593503,38031.4516148799,38509.2444533874,38993.0398499575)
for ($i = 0; $i -lt $row99.Length; $i++) {
    $ws.Cells.Item(99, 10 + $i).Value = $row99[$i]
}

$row100 = @(48406.8948342971,18934.015127415,46740.5174651586,32261.4966436238,26637.0573480246,45577.6966664176,45577.7,47921.5962564591,50386.0306195151,52977.201928004,55701.6278046219,58566.1610498258,61578.0068788139,64744.7410449404,68074.3288951533,71575.1454053834,75255.9962462813,79126.13993228909,83195.31110975859,87473.7450426932,91972.2033577031,96702.0011129302,101675.03525903,106903.8145638,112401.491075718,118181.893205549,124259.560509203,130649.780259361,137368.625897845,144432.997465451,151860.664110947,159670.308786151,167881.575239516,176515.117426427,185592.651460488,195137.010236481)
for ($i = 0; $i -lt $row100.Length; $i++) {
    $ws.Cells.Item(100, 10 + $i).Value = $row100[$i]
}

$row101 = @(972671.4157053831,984624.644550232,1007336.73713699,1050252.53767421,1025396.21661986,1035318.00405061,1035318,1073063.72562282,1112185.58862836,1152733.76036885,1194760.24134865,1238318.92791163,1283465.68135984,1330258.39959223,1378757.09135502,1429023.95319921,1481123.44924381,1535122.393847,1591090.03729143,1649098.15459337,1709221.1375497,1771536.09014077,1836122.92741143,1903064.4779571,1972446.59014631,2044358.24221585,2118891.65637985,2196142.41709897,2276209.59366159,2359195.86723395,2445207.66254236,2534355.28435629,2626753.05894734,2722519.48070562,2821777.36410136,2924654.00118683)
for ($i = 0; $i -lt $row101.Length; $i++) {
    $ws.Cells.Item(101, 10 + $i).Value = $row101[$i]
}

$row102 = @(36630.385957661,19644.0733181737,72266.6990286061,75731.84666051219,86736.58080400901,77407.4703003959,77407.47,88597.24144518111,101404.569761752,116063.283696349,132841.013516732,152044.077249449,174023.073256098,199179.28125416,227971.988648543,260926.875933719,298645.614262248,341816.850406515,391228.779670713,447783.536301983,512513.664132465,586600.9681636459,671398.091274286,768453.209986861,879538.300173489,1006681.48876013,1152204.08208762,1318762.94697192,1509399.01823208,1727592.81830846,1977328.00261571,2263164.09080495,2590319.71181993,2964767.88258707,3393345.05988227,3883875.95637994)
for ($i = 0; $i -lt $row102.Length; $i++) {
    $ws.Cells.Item(102, 10 + $i).Value = $row102[$i]
}

$row103 = @(887100.625428757,760069.013641543,828825.1093487,796529.410198479,806580.133683963,778032.335293796,778032.3,789622.522687932,801385.4030688111,813323.513196662,825439.463440991,837735.9030575569,850215.520767661,862881.045346054,875735.2462176071,888780.934062871,902020.961432656,915458.2233717639,929095.658052018,942936.247414713,956983.017822651,971239.040721872,985707.433313259,1000391.35923414,1015294.02925003,1030418.70195672,1045768.68449278,1061347.33326268,1077158.05467072,1093204.30586586,1109489.59549766,1126017.48448345,1142791.586787,1159815.57020872,1177093.15718762,1194628.12561534)
for ($i = 0; $i -lt $row103.Length; $i++) {
    $ws.Cells.Item(103, 10 + $i).Value = $row103[$i]
}

$row104 = @(381036.874589203,390797.185407436,377004.947040244,366201.743292061,390249.154506288,379341.706513841,379341.7,399355.355861375,420424.910457155,442606.071857125,465957.487228288,490540.897898451,516421.30260078,543667.129330968,572350.4162714001,602547.0022606821,634336.7273121279,667803.643711394,703036.238251387,740127.666192051,779175.997563608,820284.476464508,863561.7940396341,909122.3758605639,957086.684467681,1007581.5378741,1060740.44487348,1116703.95803836,1175620.04534227,1237644.48138823,1302941.25927801,1371683.02421125,1444051.5299607,1520238.11943062,1600444.23056902,1684881.9289711)
for ($i = 0; $i -lt $row104.Length; $i++) {
    $ws.Cells.Item(104, 10 + $i).Value = $row104[$i]
}

$row111 = @(14760.6623688479,12790.2508904471,13754.9850669999,13434.2003067164,13668.5596829101,12685.9419671159,12685.94,12636.6350133479,12587.5216547272,12538.5991793593,12489.8668453601,12441.3239137289,12392.9696483374,12344.8033159181,12296.8241860536,12249.0315311651,12201.4246265017,12154.0027501294,12106.76518292,12059.7112085402,12012.8401134408,11966.1511868459,11919.6437207422,11873.3170098679,11827.1703517024,11781.2030464556,11735.4143970571,11689.8037091456,11644.3702910587,11599.1134538222,11554.0325111394,11509.1267793813,11464.3955775757,11419.8382273969,11375.4540531559,11331.2423817896)
for ($i = 0; $i -lt $row111.Length; $i++) {
    $ws.Cells.Item(111, 10 + $i).Value = $row111[$i]
}

$row112 = @(34629.736640591,35824.4813963268,35324.7326178198,33199.8234809784,35407.0341028948,33385.1562612783,33385.16,33340.0137322168,33294.9285150769,33249.9042660225,33204.9409026071,33160.0383424962,33115.1965034662,33070.4153034049,33025.6946603111,32981.0344922946,32936.4347175756,32891.8952544854,32847.4160214652,32802.9969370669,32758.6379199524,32714.3388888935,32670.0997627719,32625.9204605792,32581.8009014162,32537.7410044933,32493.7406891302,32449.7998747555,32405.918480907,32362.0964272311,32318.3336334829,32274.6300195262,32230.9855053329,32187.4000109834,32143.8734566658,32100.4057626765)
for ($i = 0; $i -lt $row112.Length; $i++) {
    $ws.Cells.Item(112, 10 + $i).Value = $row112[$i]
}
